{"js": "// Highlight (bold + 14pt) the name, CPF and amount values that were filled\n// in from the source .xlsx, and strip the stray leading indentation spaces\n// that used to precede \"Cachoeirinha...\", the signature rule and\n// \"Assinatura\" (dropping the now pointless leading/trailing 4-space-only\n// fragments entirely).\n\nasync function highlight(para, searchText) {\n  const results = para.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].font.bold = true;\n    results.items[0].font.size = 14;\n  }\n}\n\nasync function dedent(para, oldText) {\n  const results = para.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(oldText.trim(), \"Replace\");\n    await context.sync();\n  }\n}\n\nasync function removeExact(para, literalText) {\n  const results = para.search(literalText, { matchCase: true });\n  results.load(\"items/text\");\n  await context.sync();\n  for (const item of results.items) {\n    if (item.text === literalText) {\n      item.insertText(\"\", \"Replace\");\n      await context.sync();\n    }\n  }\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Only the \"receipt\" paragraphs (every other paragraph; the ones in\n// between are blank separator paragraphs) contain the text we need to\n// touch, so detect them by looking for the constant heading text, then\n// pull out the name / CPF / amount placeholders straight from the text.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (!text || text.indexOf(\"RECIBO DE PAGAMENTO\") === -1) {\n    continue; // skip the blank separator paragraphs\n  }\n\n  const m = text.match(\n    /\\u000b {4}(.+?), inscrito\\(a\\) no CPF sob o n[\u00bao] (\\d+),[\\s\\S]*?R\\$([\\d.]+) concernente/\n  );\n  if (!m) {\n    continue;\n  }\n  const name = m[1];\n  const cpf = m[2];\n  const amount = m[3];\n\n  // 1) Bold + 14pt the \"<name>, \" fragment.\n  await highlight(para, name + \", \");\n  // 2) Bold + 14pt the CPF number.\n  await highlight(para, cpf);\n  // 3) Bold + 14pt the amount value (it only occurs once per paragraph,\n  //    right after \"R$\").\n  await highlight(para, amount);\n\n  // 4) Remove the 4-space indentation in front of the \"Cachoeirinha...\"\n  //    line, the signature rule line, and the \"Assinatura\" line.\n  await dedent(para, \"    Cachoeirinha, 29 de setembro de 2024.\");\n  await dedent(para, \"    _________________________________________________\");\n  await dedent(para, \"    Assinatura\");\n\n  // 5) Two spaces-only (\"    \") fragments are left over: the indent that\n  //    used to sit right before the name, and the paragraph's old\n  //    trailing spaces-only line. Both are now pointless, drop them.\n  await removeExact(para, \"    \");\n}\n", "ps1": "# Highlight (bold + 14pt) the name, CPF and amount values that were filled\n# in from the source .xlsx, and strip the stray leading indentation spaces\n# that used to precede \"Cachoeirinha...\", the signature rule and\n# \"Assinatura\" (dropping the now pointless leading/trailing 4-space-only\n# runs entirely).\n\n$d = $word.ActiveDocument\n\nfunction BoldRange($para, $searchText) {\n    $r = $para.Range.Duplicate\n    $r.Find.ClearFormatting()\n    $r.Find.Text = $searchText\n    $r.Find.MatchCase = $true\n    $r.Find.Forward = $true\n    $r.Find.Wrap = 1\n    $ok = $r.Find.Execute()\n    if ($ok) {\n        $r.Font.Bold = 1\n        $r.Font.Size = 14\n    }\n    return $ok\n}\n\nfunction Dedent($para, $oldText) {\n    $r = $para.Range.Duplicate\n    $r.Find.ClearFormatting()\n    $r.Find.Text = $oldText\n    $r.Find.Replacement.Text = $oldText.Trim()\n    $r.Find.MatchCase = $true\n    $r.Find.Forward = $true\n    $r.Find.Wrap = 1\n    $r.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n\nfunction RemoveExact($para, $literalText) {\n    $iterations = 0\n    while ($iterations -lt 10) {\n        $r = $para.Range.Duplicate\n        $r.Find.ClearFormatting()\n        $r.Find.Text = $literalText\n        $r.Find.MatchCase = $true\n        $r.Find.Forward = $true\n        $r.Find.Wrap = 1\n        $ok = $r.Find.Execute()\n        if (-not $ok) { break }\n        if ($r.Text -ne $literalText) { break }\n        $r.Text = \"\"\n        $iterations = $iterations + 1\n    }\n}\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -notmatch 'RECIBO DE PAGAMENTO') {\n        continue\n    }\n    if ($text -notmatch '    (.+?), inscrito\\(a\\) no CPF sob o n. (\\d+),[\\s\\S]*?R\\$([\\d.]+) concernente') {\n        continue\n    }\n\n    $name = $matches[1]\n    $cpf = $matches[2]\n    $amount = $matches[3]\n\n    BoldRange $p \"$name, \" | Out-Null\n    BoldRange $p $cpf | Out-Null\n    BoldRange $p $amount | Out-Null\n\n    Dedent $p \"    Cachoeirinha, 29 de setembro de 2024.\"\n    Dedent $p \"    _________________________________________________\"\n    Dedent $p \"    Assinatura\"\n\n    # Remove the leftover spaces-only runs: the 4-space indent that used to\n    # sit right before the name, and the trailing spaces-only line.\n    RemoveExact $p \"    \"\n}\n\nWrite-Output \"done\"\n"}
